# Update column F ("dSF") values per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 2
    3  = -2
    5  = 5
    6  = 3
    7  = -1
    8  = -3
    9  = -2
    10 = 1
    11 = -1
    12 = -1
    13 = -1
    14 = 4
    15 = 8
    16 = -1
    17 = 1
    18 = -1
    19 = -3
    20 = 4
    22 = 6
    23 = 1
    24 = 2
    25 = 2
    26 = 1
    27 = 1
    29 = -2
    30 = -6
    32 = 2
    33 = 4
    34 = -2
    35 = 3
    36 = 2
    37 = 5
    38 = 7
    39 = -5
    40 = 1
    41 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
